# The "1_タイトルのみ" (Title Only) slide layout (ppt/slideLayouts/slideLayout4.xml,
# the 4th CustomLayout of the slide master) carries an author-only helper
# group "Group 46" (a red content-area rectangle + a red "Note:" instructions
# text box). It was userDrawn scaffolding, not meant to ship — remove it.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$layout = $master.CustomLayouts.Item(4)

$helper = $layout.Shapes.Item("Group 46")
$helper.Delete()
